$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.628848666666667
$ws.Cells.Item(2, 8).Value = 7.886546000000001
$ws.Cells.Item(2, 9).Value = 0.04622248078033103
$ws.Cells.Item(2, 10).Value = 0.04850184447997802
$ws.Cells.Item(2, 13).Value = 63.46725166666666
$ws.Cells.Item(2, 14).Value = 190.401755
$ws.Cells.Item(2, 15).Value = 0.2354497988808272
$ws.Cells.Item(2, 16).Value = 0.2397164477183668
$ws.Cells.Item(2, 17).Value = 166.8457999209144
$ws.Cells.Item(2, 18).Value = 1501.61219928823
$ws.Cells.Item(2, 19).Value = 0.01088307380350184
$ws.Cells.Item(2, 20).Value = 0.01162668986652901

$ws.Cells.Item(3, 7).Value = 2.628848666666667
$ws.Cells.Item(3, 8).Value = 7.886546000000001
$ws.Cells.Item(3, 9).Value = 0.04622248078033103
$ws.Cells.Item(3, 10).Value = 0.04850184447997802
$ws.Cells.Item(3, 13).Value = 47.980825
$ws.Cells.Item(3, 14).Value = 143.942475
$ws.Cells.Item(3, 15).Value = 0.1779985000094065
$ws.Cells.Item(3, 16).Value = 0.1812240584798697
$ws.Cells.Item(3, 17).Value = 126.1343278268167
$ws.Cells.Item(3, 18).Value = 1135.20895044135
$ws.Cells.Item(3, 19).Value = 0.008227532245612543
$ws.Cells.Item(3, 20).Value = 0.008789701100421083

$ws.Cells.Item(4, 7).Value = 2.628848666666667
$ws.Cells.Item(4, 8).Value = 7.886546000000001
$ws.Cells.Item(4, 9).Value = 0.04622248078033103
$ws.Cells.Item(4, 10).Value = 0.04850184447997802
$ws.Cells.Item(4, 13).Value = 64.53809233333334
$ws.Cells.Item(4, 14).Value = 193.614277
$ws.Cells.Item(4, 15).Value = 0.2394223865221556
$ws.Cells.Item(4, 16).Value = 0.243761023683841
$ws.Cells.Item(4, 17).Value = 169.6608779796936
$ws.Cells.Item(4, 18).Value = 1526.947901817242
$ws.Cells.Item(4, 19).Value = 0.01106669665940133
$ws.Cells.Item(4, 20).Value = 0.0118228592609939

$ws.Cells.Item(5, 7).Value = 2.628848666666667
$ws.Cells.Item(5, 8).Value = 7.886546000000001
$ws.Cells.Item(5, 9).Value = 0.04622248078033103
$ws.Cells.Item(5, 10).Value = 0.04850184447997802
$ws.Cells.Item(5, 13).Value = 14.3933435
$ws.Cells.Item(5, 14).Value = 28.786687
$ws.Cells.Item(5, 15).Value = 0.0533961963580272
$ws.Cells.Item(5, 16).Value = 0.03624253541791403
$ws.Cells.Item(5, 17).Value = 37.83792186885034
$ws.Cells.Item(5, 18).Value = 227.027531213102
$ws.Cells.Item(5, 19).Value = 0.002468104659901694
$ws.Cells.Item(5, 20).Value = 0.001757829816399762

$ws.Cells.Item(6, 7).Value = 2.628848666666667
$ws.Cells.Item(6, 8).Value = 7.886546000000001
$ws.Cells.Item(6, 9).Value = 0.04622248078033103
$ws.Cells.Item(6, 10).Value = 0.04850184447997802
$ws.Cells.Item(6, 13).Value = 79.17795566666666
$ws.Cells.Item(6, 14).Value = 237.533867
$ws.Cells.Item(6, 15).Value = 0.2937331182295834
$ws.Cells.Item(6, 16).Value = 0.2990559347000084
$ws.Cells.Item(6, 17).Value = 208.1468631837091
$ws.Cells.Item(6, 18).Value = 1873.321768653382
$ws.Cells.Item(6, 19).Value = 0.01357707341191362
$ws.Cells.Item(6, 20).Value = 0.01450476443563427

$ws.Cells.Item(7, 7).Value = 26.85202466666667
$ws.Cells.Item(7, 8).Value = 80.556074
$ws.Cells.Item(7, 9).Value = 0.472133375270229
$ws.Cells.Item(7, 10).Value = 0.4954156322762335
$ws.Cells.Item(7, 13).Value = 63.46725166666666
$ws.Cells.Item(7, 14).Value = 190.401755
$ws.Cells.Item(7, 15).Value = 0.2354497988808272
$ws.Cells.Item(7, 16).Value = 0.2397164477183668
$ws.Cells.Item(7, 17).Value = 1704.224207278874
$ws.Cells.Item(7, 18).Value = 15338.01786550987
$ws.Cells.Item(7, 19).Value = 0.1111637082523016
$ws.Cells.Item(7, 20).Value = 0.1187592755134074

$ws.Cells.Item(8, 7).Value = 26.85202466666667
$ws.Cells.Item(8, 8).Value = 80.556074
$ws.Cells.Item(8, 9).Value = 0.472133375270229
$ws.Cells.Item(8, 10).Value = 0.4954156322762335
$ws.Cells.Item(8, 13).Value = 47.980825
$ws.Cells.Item(8, 14).Value = 143.942475
$ws.Cells.Item(8, 15).Value = 0.1779985000094065
$ws.Cells.Item(8, 16).Value = 0.1812240584798697
$ws.Cells.Item(8, 17).Value = 1288.382296427017
$ws.Cells.Item(8, 18).Value = 11595.44066784315
$ws.Cells.Item(8, 19).Value = 0.08403903260247896
$ws.Cells.Item(8, 20).Value = 0.08978123151546977

$ws.Cells.Item(9, 7).Value = 26.85202466666667
$ws.Cells.Item(9, 8).Value = 80.556074
$ws.Cells.Item(9, 9).Value = 0.472133375270229
$ws.Cells.Item(9, 10).Value = 0.4954156322762335
$ws.Cells.Item(9, 13).Value = 64.53809233333334
$ws.Cells.Item(9, 14).Value = 193.614277
$ws.Cells.Item(9, 15).Value = 0.2394223865221556
$ws.Cells.Item(9, 16).Value = 0.243761023683841
$ws.Cells.Item(9, 17).Value = 1732.978447274278
$ws.Cells.Item(9, 18).Value = 15596.8060254685
$ws.Cells.Item(9, 19).Value = 0.1130392994639587
$ws.Cells.Item(9, 20).Value = 0.120763021672632

$ws.Cells.Item(10, 7).Value = 26.85202466666667
$ws.Cells.Item(10, 8).Value = 80.556074
$ws.Cells.Item(10, 9).Value = 0.472133375270229
$ws.Cells.Item(10, 10).Value = 0.4954156322762335
$ws.Cells.Item(10, 13).Value = 14.3933435
$ws.Cells.Item(10, 14).Value = 28.786687
$ws.Cells.Item(10, 15).Value = 0.0533961963580272
$ws.Cells.Item(10, 16).Value = 0.03624253541791403
$ws.Cells.Item(10, 17).Value = 386.4904146978063
$ws.Cells.Item(10, 18).Value = 2318.942488186838
$ws.Cells.Item(10, 19).Value = 0.02521012641310729
$ws.Cells.Item(10, 20).Value = 0.01795511859935967

$ws.Cells.Item(11, 7).Value = 26.85202466666667
$ws.Cells.Item(11, 8).Value = 80.556074
$ws.Cells.Item(11, 9).Value = 0.472133375270229
$ws.Cells.Item(11, 10).Value = 0.4954156322762335
$ws.Cells.Item(11, 13).Value = 79.17795566666666
$ws.Cells.Item(11, 14).Value = 237.533867
$ws.Cells.Item(11, 15).Value = 0.2937331182295834
$ws.Cells.Item(11, 16).Value = 0.2990559347000084
$ws.Cells.Item(11, 17).Value = 2126.088418617573
$ws.Cells.Item(11, 18).Value = 19134.79576755816
$ws.Cells.Item(11, 19).Value = 0.1386812085383824
$ws.Cells.Item(11, 20).Value = 0.1481569849753646

$ws.Cells.Item(12, 7).Value = 11.96574466666667
$ws.Cells.Item(12, 8).Value = 35.897234
$ws.Cells.Item(12, 9).Value = 0.2103911152781009
$ws.Cells.Item(12, 10).Value = 0.2207661073338543
$ws.Cells.Item(12, 13).Value = 63.46725166666666
$ws.Cells.Item(12, 14).Value = 190.401755
$ws.Cells.Item(12, 15).Value = 0.2354497988808272
$ws.Cells.Item(12, 16).Value = 0.2397164477183668
$ws.Cells.Item(12, 17).Value = 759.4329281384076
$ws.Cells.Item(12, 18).Value = 6834.896353245669
$ws.Cells.Item(12, 19).Value = 0.04953654577854179
$ws.Cells.Item(12, 20).Value = 0.05292126702668323

$ws.Cells.Item(13, 7).Value = 11.96574466666667
$ws.Cells.Item(13, 8).Value = 35.897234
$ws.Cells.Item(13, 9).Value = 0.2103911152781009
$ws.Cells.Item(13, 10).Value = 0.2207661073338543
$ws.Cells.Item(13, 13).Value = 47.980825
$ws.Cells.Item(13, 14).Value = 143.942475
$ws.Cells.Item(13, 15).Value = 0.1779985000094065
$ws.Cells.Item(13, 16).Value = 0.1812240584798697
$ws.Cells.Item(13, 17).Value = 574.1263008460166
$ws.Cells.Item(13, 18).Value = 5167.13670761415
$ws.Cells.Item(13, 19).Value = 0.03744930293480807
$ws.Cells.Item(13, 20).Value = 0.0400081299458436

$ws.Cells.Item(14, 7).Value = 11.96574466666667
$ws.Cells.Item(14, 8).Value = 35.897234
$ws.Cells.Item(14, 9).Value = 0.2103911152781009
$ws.Cells.Item(14, 10).Value = 0.2207661073338543
$ws.Cells.Item(14, 13).Value = 64.53809233333334
$ws.Cells.Item(14, 14).Value = 193.614277
$ws.Cells.Item(14, 15).Value = 0.2394223865221556
$ws.Cells.Item(14, 16).Value = 0.243761023683841
$ws.Cells.Item(14, 17).Value = 772.2463341344243
$ws.Cells.Item(14, 18).Value = 6950.217007209818
$ws.Cells.Item(14, 19).Value = 0.05037234292294088
$ws.Cells.Item(14, 20).Value = 0.05381417231839705

$ws.Cells.Item(15, 7).Value = 11.96574466666667
$ws.Cells.Item(15, 8).Value = 35.897234
$ws.Cells.Item(15, 9).Value = 0.2103911152781009
$ws.Cells.Item(15, 10).Value = 0.2207661073338543
$ws.Cells.Item(15, 13).Value = 14.3933435
$ws.Cells.Item(15, 14).Value = 28.786687
$ws.Cells.Item(15, 15).Value = 0.0533961963580272
$ws.Cells.Item(15, 16).Value = 0.03624253541791403
$ws.Cells.Item(15, 17).Value = 172.2270732206263
$ws.Cells.Item(15, 18).Value = 1033.362439323758
$ws.Cells.Item(15, 19).Value = 0.01123408530337381
$ws.Cells.Item(15, 20).Value = 0.008001123464122224

$ws.Cells.Item(16, 7).Value = 11.96574466666667
$ws.Cells.Item(16, 8).Value = 35.897234
$ws.Cells.Item(16, 9).Value = 0.2103911152781009
$ws.Cells.Item(16, 10).Value = 0.2207661073338543
$ws.Cells.Item(16, 13).Value = 79.17795566666666
$ws.Cells.Item(16, 14).Value = 237.533867
$ws.Cells.Item(16, 15).Value = 0.2937331182295834
$ws.Cells.Item(16, 16).Value = 0.2990559347000084
$ws.Cells.Item(16, 17).Value = 947.4232007359864
$ws.Cells.Item(16, 18).Value = 8526.808806623876
$ws.Cells.Item(16, 19).Value = 0.06179883833843632
$ws.Cells.Item(16, 20).Value = 0.06602141457880817

$ws.Cells.Item(17, 7).Value = 8.018423
$ws.Cells.Item(17, 8).Value = 16.036846
$ws.Cells.Item(17, 9).Value = 0.1409862072722574
$ws.Cells.Item(17, 10).Value = 0.09862576223372788
$ws.Cells.Item(17, 13).Value = 63.46725166666666
$ws.Cells.Item(17, 14).Value = 190.401755
$ws.Cells.Item(17, 15).Value = 0.2354497988808272
$ws.Cells.Item(17, 16).Value = 0.2397164477183668
$ws.Cells.Item(17, 17).Value = 508.9072705107883
$ws.Cells.Item(17, 18).Value = 3053.44362306473
$ws.Cells.Item(17, 19).Value = 0.03319517414722364
$ws.Cells.Item(17, 20).Value = 0.0236422173761855

$ws.Cells.Item(18, 7).Value = 8.018423
$ws.Cells.Item(18, 8).Value = 16.036846
$ws.Cells.Item(18, 9).Value = 0.1409862072722574
$ws.Cells.Item(18, 10).Value = 0.09862576223372788
$ws.Cells.Item(18, 13).Value = 47.980825
$ws.Cells.Item(18, 14).Value = 143.942475
$ws.Cells.Item(18, 15).Value = 0.1779985000094065
$ws.Cells.Item(18, 16).Value = 0.1812240584798697
$ws.Cells.Item(18, 17).Value = 384.730550738975
$ws.Cells.Item(18, 18).Value = 2308.38330443385
$ws.Cells.Item(18, 19).Value = 0.0250953334164771
$ws.Cells.Item(18, 20).Value = 0.01787336090266682

$ws.Cells.Item(19, 7).Value = 8.018423
$ws.Cells.Item(19, 8).Value = 16.036846
$ws.Cells.Item(19, 9).Value = 0.1409862072722574
$ws.Cells.Item(19, 10).Value = 0.09862576223372788
$ws.Cells.Item(19, 13).Value = 64.53809233333334
$ws.Cells.Item(19, 14).Value = 193.614277
$ws.Cells.Item(19, 15).Value = 0.2394223865221556
$ws.Cells.Item(19, 16).Value = 0.243761023683841
$ws.Cells.Item(19, 17).Value = 517.4937239417237
$ws.Cells.Item(19, 18).Value = 3104.962343650343
$ws.Cells.Item(19, 19).Value = 0.03375525421183117
$ws.Cells.Item(19, 20).Value = 0.02404111676369261

$ws.Cells.Item(20, 7).Value = 8.018423
$ws.Cells.Item(20, 8).Value = 16.036846
$ws.Cells.Item(20, 9).Value = 0.1409862072722574
$ws.Cells.Item(20, 10).Value = 0.09862576223372788
$ws.Cells.Item(20, 13).Value = 14.3933435
$ws.Cells.Item(20, 14).Value = 28.786687
$ws.Cells.Item(20, 15).Value = 0.0533961963580272
$ws.Cells.Item(20, 16).Value = 0.03624253541791403
$ws.Cells.Item(20, 17).Value = 115.4119165673005
$ws.Cells.Item(20, 18).Value = 461.647666269202
$ws.Cells.Item(20, 19).Value = 0.007528127207282981
$ws.Cells.Item(20, 20).Value = 0.00357444768087465

$ws.Cells.Item(21, 7).Value = 8.018423
$ws.Cells.Item(21, 8).Value = 16.036846
$ws.Cells.Item(21, 9).Value = 0.1409862072722574
$ws.Cells.Item(21, 10).Value = 0.09862576223372788
$ws.Cells.Item(21, 13).Value = 79.17795566666666
$ws.Cells.Item(21, 14).Value = 237.533867
$ws.Cells.Item(21, 15).Value = 0.2937331182295834
$ws.Cells.Item(21, 16).Value = 0.2990559347000084
$ws.Cells.Item(21, 17).Value = 634.8823408105803
$ws.Cells.Item(21, 18).Value = 3809.294044863482
$ws.Cells.Item(21, 19).Value = 0.04141231828944254
$ws.Cells.Item(21, 20).Value = 0.02949461951030827

$ws.Cells.Item(22, 7).Value = 7.408770666666666
$ws.Cells.Item(22, 8).Value = 22.226312
$ws.Cells.Item(22, 9).Value = 0.1302668213990815
$ws.Cells.Item(22, 10).Value = 0.1366906536762062
$ws.Cells.Item(22, 13).Value = 63.46725166666666
$ws.Cells.Item(22, 14).Value = 190.401755
$ws.Cells.Item(22, 15).Value = 0.2354497988808272
$ws.Cells.Item(22, 16).Value = 0.2397164477183668
$ws.Cells.Item(22, 17).Value = 470.2143124419511
$ws.Cells.Item(22, 18).Value = 4231.928811977559
$ws.Cells.Item(22, 19).Value = 0.03067129689925839
$ws.Cells.Item(22, 20).Value = 0.03276699793556166

$ws.Cells.Item(23, 7).Value = 7.408770666666666
$ws.Cells.Item(23, 8).Value = 22.226312
$ws.Cells.Item(23, 9).Value = 0.1302668213990815
$ws.Cells.Item(23, 10).Value = 0.1366906536762062
$ws.Cells.Item(23, 13).Value = 47.980825
$ws.Cells.Item(23, 14).Value = 143.942475
$ws.Cells.Item(23, 15).Value = 0.1779985000094065
$ws.Cells.Item(23, 16).Value = 0.1812240584798697
$ws.Cells.Item(23, 17).Value = 355.4789288224667
$ws.Cells.Item(23, 18).Value = 3199.3103594022
$ws.Cells.Item(23, 19).Value = 0.02318729881002976
$ws.Cells.Item(23, 20).Value = 0.02477163501546841

$ws.Cells.Item(24, 7).Value = 7.408770666666666
$ws.Cells.Item(24, 8).Value = 22.226312
$ws.Cells.Item(24, 9).Value = 0.1302668213990815
$ws.Cells.Item(24, 10).Value = 0.1366906536762062
$ws.Cells.Item(24, 13).Value = 64.53809233333334
$ws.Cells.Item(24, 14).Value = 193.614277
$ws.Cells.Item(24, 15).Value = 0.2394223865221556
$ws.Cells.Item(24, 16).Value = 0.243761023683841
$ws.Cells.Item(24, 17).Value = 478.1479253618249
$ws.Cells.Item(24, 18).Value = 4303.331328256424
$ws.Cells.Item(24, 19).Value = 0.03118879326402351
$ws.Cells.Item(24, 20).Value = 0.03331985366812541

$ws.Cells.Item(25, 7).Value = 7.408770666666666
$ws.Cells.Item(25, 8).Value = 22.226312
$ws.Cells.Item(25, 9).Value = 0.1302668213990815
$ws.Cells.Item(25, 10).Value = 0.1366906536762062
$ws.Cells.Item(25, 13).Value = 14.3933435
$ws.Cells.Item(25, 14).Value = 28.786687
$ws.Cells.Item(25, 15).Value = 0.0533961963580272
$ws.Cells.Item(25, 16).Value = 0.03624253541791403
$ws.Cells.Item(25, 17).Value = 106.6369811180573
$ws.Cells.Item(25, 18).Value = 639.821886708344
$ws.Cells.Item(25, 19).Value = 0.006955752774361416
$ws.Cells.Item(25, 20).Value = 0.004954015857157723

$ws.Cells.Item(26, 7).Value = 7.408770666666666
$ws.Cells.Item(26, 8).Value = 22.226312
$ws.Cells.Item(26, 9).Value = 0.1302668213990815
$ws.Cells.Item(26, 10).Value = 0.1366906536762062
$ws.Cells.Item(26, 13).Value = 79.17795566666666
$ws.Cells.Item(26, 14).Value = 237.533867
$ws.Cells.Item(26, 15).Value = 0.2937331182295834
$ws.Cells.Item(26, 16).Value = 0.2990559347000084
$ws.Cells.Item(26, 17).Value = 586.6113153898337
$ws.Cells.Item(26, 18).Value = 5279.501838508504
$ws.Cells.Item(26, 19).Value = 0.03826367965140844
$ws.Cells.Item(26, 20).Value = 0.04087815119989298
